$d = $word.ActiveDocument

$replacements = @(
    @("217÷3=72, 1", "252÷4=63, 0"),
    @("762÷2=381, 0", "710÷4=177, 2"),
    @("207÷7=29, 4", "662÷3=220, 2"),
    @("405÷7=57, 6", "355÷2=177, 1"),
    @("776÷9=86, 2", "303÷4=75, 3"),
    @("862÷8=107, 6", "581÷4=145, 1"),
    @("461÷4=115, 1", "395÷5=79, 0"),
    @("437÷2=218, 1", "770÷2=385, 0"),
    @("335÷8=41, 7", "335÷4=83, 3"),
    @("310÷9=34, 4", "842÷9=93, 5"),
    @("735÷9=81, 6", "592÷7=84, 4"),
    @("598÷3=199, 1", "925÷7=132, 1"),
    @("598÷9=66, 4", "157÷5=31, 2"),
    @("218÷5=43, 3", "841÷2=420, 1"),
    @("843÷8=105, 3", "464÷9=51, 5"),
    @("296÷6=49, 2", "185÷2=92, 1"),
    @("178÷3=59, 1", "180÷9=20, 0"),
    @("810÷5=162, 0", "977÷5=195, 2"),
    @("144÷2=72, 0", "239÷8=29, 7"),
    @("437÷5=87, 2", "653÷6=108, 5"),
    @("873÷3=291, 0", "294÷4=73, 2"),
    @("345÷9=38, 3", "136÷2=68, 0"),
    @("967÷9=107, 4", "251÷9=27, 8"),
    @("441÷3=147, 0", "896÷8=112, 0"),
    @("180÷7=25, 5", "152÷9=16, 8")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
